# Fix the algorithm/conditions on filtering the status of candidates
# Re-align candidate rows (D/E/F columns) for the "groundcover" SDR role (rows 4-6)
# and the "Orca Ai" Singapore AE role (rows 10-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 4-6 (groundcover / SDR candidates) ---
$ws.Range("D4").Value = "Cesar Castillo"
$ws.Range("E4").Value = "CV Sent"
$ws.Range("F4").Value = 45986

$ws.Range("D5").Value = "JJ Valderrama"
$ws.Range("E5").Value = "CV Sent"
$ws.Range("F5").Value = 45986

$ws.Range("D6").Value = "Tyler Drago"
$ws.Range("E6").Value = "3rd Interview"
$ws.Range("F6").Value = 45992

# --- Rows 10-11 (Orca Ai Singapore AE candidates) ---
$ws.Range("D10").Value = "Gabriel Wong"
$ws.Range("E10").Value = "2nd Interview"
$ws.Range("F10").Value = 45992

$ws.Range("D11").Value = "Jodie Yao"
$ws.Range("E11").Value = "2nd Interview"
$ws.Range("F11").Value = 45995
